$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 11494836
$ws.Cells.Item(19, 9).Value = 33333620
$ws.Cells.Item(19, 10).Value = 739.0526
$ws.Cells.Item(19, 11).Value = 33333620
$ws.Cells.Item(19, 12).Value = 739.0526
$ws.Cells.Item(19, 13).Value = -33333445
$ws.Cells.Item(19, 14).Value = -1089.0526

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(63, 8).Value = 38271
$ws.Cells.Item(63, 10).Value = 38271
$ws.Cells.Item(63, 12).Value = 38271
$ws.Cells.Item(63, 14).Value = -39519

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(66, 8).Value = 38271
$ws.Cells.Item(66, 10).Value = 38271
$ws.Cells.Item(66, 12).Value = 114813
$ws.Cells.Item(66, 14).Value = -121053

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 1854.9429
$ws.Cells.Item(132, 9).Value = 1924.9395
$ws.Cells.Item(132, 10).Value = 700
$ws.Cells.Item(132, 11).Value = 5774.818499999999
$ws.Cells.Item(132, 12).Value = 2100
$ws.Cells.Item(132, 13).Value = -3244.818499999999
$ws.Cells.Item(132, 14).Value = -7160

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1891.8182
$ws.Cells.Item(2, 9).Value = 1896
$ws.Cells.Item(2, 10).Value = 1850
$ws.Cells.Item(2, 11).Value = 1896
$ws.Cells.Item(2, 12).Value = 1850
$ws.Cells.Item(2, 13).Value = -1783
$ws.Cells.Item(2, 14).Value = -2076

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(70, 8).Value = 59800
$ws.Cells.Item(70, 10).Value = 59800
$ws.Cells.Item(70, 12).Value = 59800
$ws.Cells.Item(70, 14).Value = -60340

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(73, 8).Value = 59800
$ws.Cells.Item(73, 10).Value = 59800
$ws.Cells.Item(73, 12).Value = 59800
$ws.Cells.Item(73, 14).Value = -61672

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(101, 8).Value = 33900
$ws.Cells.Item(101, 10).Value = 33900
$ws.Cells.Item(101, 12).Value = 33900
$ws.Cells.Item(101, 14).Value = -40390

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 1891.8182
$ws.Cells.Item(116, 9).Value = 1896
$ws.Cells.Item(116, 10).Value = 1850
$ws.Cells.Item(116, 11).Value = 1896
$ws.Cells.Item(116, 12).Value = 1850
$ws.Cells.Item(116, 13).Value = 398
$ws.Cells.Item(116, 14).Value = -6438

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 2383.1667
$ws.Cells.Item(132, 9).Value = 1484.238
$ws.Cells.Item(132, 10).Value = 4480.6665
$ws.Cells.Item(132, 11).Value = 4452.714
$ws.Cells.Item(132, 12).Value = 13441.9995
$ws.Cells.Item(132, 13).Value = -1922.714
$ws.Cells.Item(132, 14).Value = -18501.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1891.8182
$ws.Cells.Item(3, 9).Value = 1896
$ws.Cells.Item(3, 10).Value = 1850
$ws.Cells.Item(3, 11).Value = 1896
$ws.Cells.Item(3, 12).Value = 1850
$ws.Cells.Item(3, 13).Value = -1782
$ws.Cells.Item(3, 14).Value = -2078

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(64, 8).Value = 669
$ws.Cells.Item(64, 10).Value = 1007
$ws.Cells.Item(64, 12).Value = 1007
$ws.Cells.Item(64, 14).Value = -1457

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(67, 8).Value = 669
$ws.Cells.Item(67, 10).Value = 1007
$ws.Cells.Item(67, 12).Value = 1007
$ws.Cells.Item(67, 14).Value = -2567

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(102, 8).Value = 4999.5
$ws.Cells.Item(102, 9).Value = 4999.5
$ws.Cells.Item(102, 11).Value = 4999.5
$ws.Cells.Item(102, 13).Value = -1754.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1456.3572
$ws.Cells.Item(107, 9).Value = 1008.25
$ws.Cells.Item(107, 10).Value = 1635.6
$ws.Cells.Item(107, 11).Value = 1008.25
$ws.Cells.Item(107, 12).Value = 1635.6
$ws.Cells.Item(107, 13).Value = 911.75
$ws.Cells.Item(107, 14).Value = -5475.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4244.787
$ws.Cells.Item(31, 9).Value = 5052.846
$ws.Cells.Item(31, 10).Value = 3244.3333
$ws.Cells.Item(31, 11).Value = 5052.846
$ws.Cells.Item(31, 12).Value = 3244.3333
$ws.Cells.Item(31, 13).Value = -4757.846
$ws.Cells.Item(31, 14).Value = -3834.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 4244.787
$ws.Cells.Item(34, 9).Value = 5052.846
$ws.Cells.Item(34, 10).Value = 3244.3333
$ws.Cells.Item(34, 11).Value = 5052.846
$ws.Cells.Item(34, 12).Value = 3244.3333
$ws.Cells.Item(34, 13).Value = -4850.846
$ws.Cells.Item(34, 14).Value = -3648.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(60, 8).Value = 30000
$ws.Cells.Item(60, 10).Value = 30000
$ws.Cells.Item(60, 12).Value = 30000
$ws.Cells.Item(60, 14).Value = -31022

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 2271.913
$ws.Cells.Item(99, 9).Value = 1777.8462
$ws.Cells.Item(99, 10).Value = 2914.2
$ws.Cells.Item(99, 11).Value = 1777.8462
$ws.Cells.Item(99, 12).Value = 2914.2
$ws.Cells.Item(99, 13).Value = -279.8462
$ws.Cells.Item(99, 14).Value = -5910.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(103, 8).Value = 15821.223
$ws.Cells.Item(103, 9).Value = 11548.875
$ws.Cells.Item(103, 10).Value = 50000
$ws.Cells.Item(103, 11).Value = 11548.875
$ws.Cells.Item(103, 12).Value = 50000
$ws.Cells.Item(103, 13).Value = -10376.875
$ws.Cells.Item(103, 14).Value = -52344

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 2271.913
$ws.Cells.Item(126, 9).Value = 1777.8462
$ws.Cells.Item(126, 10).Value = 2914.2
$ws.Cells.Item(126, 11).Value = 5333.5386
$ws.Cells.Item(126, 12).Value = 8742.599999999999
$ws.Cells.Item(126, 13).Value = -2863.5386
$ws.Cells.Item(126, 14).Value = -13682.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(131, 8).Value = 29826
$ws.Cells.Item(131, 10).Value = 29826
$ws.Cells.Item(131, 12).Value = 29826
$ws.Cells.Item(131, 14).Value = -39906

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 5208746
$ws.Cells.Item(5, 9).Value = 407.56
$ws.Cells.Item(5, 10).Value = 23809956
$ws.Cells.Item(5, 11).Value = 1222.68
$ws.Cells.Item(5, 12).Value = 71429868
$ws.Cells.Item(5, 13).Value = -1110.68
$ws.Cells.Item(5, 14).Value = -71430092

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value = 725
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 725
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 12).Value = 2175
$ws.Cells.Item(92, 13).Value = $null
$ws.Cells.Item(92, 14).Value = -4671

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 1130.8235
$ws.Cells.Item(107, 9).Value = 405.8
$ws.Cells.Item(107, 11).Value = 1217.4
$ws.Cells.Item(107, 13).Value = 702.5999999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(135, 8).Value = 5208746
$ws.Cells.Item(135, 9).Value = 407.56
$ws.Cells.Item(135, 10).Value = 23809956
$ws.Cells.Item(135, 11).Value = 3668.04
$ws.Cells.Item(135, 12).Value = 214289604
$ws.Cells.Item(135, 13).Value = -1133.04
$ws.Cells.Item(135, 14).Value = -214294674

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5442.149
$ws.Cells.Item(70, 9).Value = 4949.05
$ws.Cells.Item(70, 10).Value = 5807.407
$ws.Cells.Item(70, 11).Value = 4949.05
$ws.Cells.Item(70, 12).Value = 5807.407
$ws.Cells.Item(70, 13).Value = -4679.05
$ws.Cells.Item(70, 14).Value = -6347.407

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 5442.149
$ws.Cells.Item(73, 9).Value = 4949.05
$ws.Cells.Item(73, 10).Value = 5807.407
$ws.Cells.Item(73, 11).Value = 4949.05
$ws.Cells.Item(73, 12).Value = 5807.407
$ws.Cells.Item(73, 13).Value = -4013.05
$ws.Cells.Item(73, 14).Value = -7679.407

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2732.4443
$ws.Cells.Item(132, 9).Value = 2449.125
$ws.Cells.Item(132, 11).Value = 7347.375
$ws.Cells.Item(132, 13).Value = -4817.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 445537.4
$ws.Cells.Item(61, 9).Value = 11243.409
$ws.Cells.Item(61, 11).Value = 11243.409
$ws.Cells.Item(61, 13).Value = -11041.409

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 445537.4
$ws.Cells.Item(113, 9).Value = 11243.409
$ws.Cells.Item(113, 11).Value = 11243.409
$ws.Cells.Item(113, 13).Value = -9073.409

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 5978.183
$ws.Cells.Item(122, 9).Value = 5221.095
$ws.Cells.Item(122, 11).Value = 15663.285
$ws.Cells.Item(122, 13).Value = -13213.285

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 6264.625
$ws.Cells.Item(132, 9).Value = 6881.087
$ws.Cells.Item(132, 11).Value = 20643.261
$ws.Cells.Item(132, 13).Value = -18113.261

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 786.75
$ws.Cells.Item(96, 9).Value = 800.25
$ws.Cells.Item(96, 10).Value = 773.25
$ws.Cells.Item(96, 11).Value = 773.25
$ws.Cells.Item(96, 12).Value = 773.25
$ws.Cells.Item(96, 13).Value = 572.75
$ws.Cells.Item(96, 14).Value = -3519.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2002.6364
$ws.Cells.Item(132, 9).Value = 920.4194
$ws.Cells.Item(132, 11).Value = 2761.2582
$ws.Cells.Item(132, 13).Value = -231.2582000000002
